$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3874.1667
$ws.Range("J17").Value = 3874.1667
$ws.Range("L17").Value = 11622.5001
$ws.Range("N17").Value = -11958.5001
$ws.Range("H40").Value = 3292.1875
$ws.Range("I40").Value = 1898.0769
$ws.Range("J40").Value = 9333.333000000001
$ws.Range("K40").Value = 1898.0769
$ws.Range("L40").Value = 9333.333000000001
$ws.Range("M40").Value = -1723.0769
$ws.Range("N40").Value = -9683.333000000001
$ws.Range("H70").Value = 1381.3462
$ws.Range("I70").Value = 1195.5
$ws.Range("J70").Value = 1497.5
$ws.Range("K70").Value = 3586.5
$ws.Range("L70").Value = 4492.5
$ws.Range("M70").Value = -3316.5
$ws.Range("N70").Value = -5032.5
$ws.Range("H73").Value = 1381.3462
$ws.Range("I73").Value = 1195.5
$ws.Range("J73").Value = 1497.5
$ws.Range("K73").Value = 3586.5
$ws.Range("L73").Value = 4492.5
$ws.Range("M73").Value = -2650.5
$ws.Range("N73").Value = -6364.5
$ws.Range("H92").Value = 924.53845
$ws.Range("I92").Value = 889.9091
$ws.Range("J92").Value = 1115
$ws.Range("K92").Value = 889.9091
$ws.Range("L92").Value = 1115
$ws.Range("M92").Value = 358.0909
$ws.Range("N92").Value = -3611
$ws.Range("H104").Value = 1291
$ws.Range("I104").Value = 1291
$ws.Range("K104").Value = 3873
$ws.Range("M104").Value = -2126
$ws.Range("H107").Value = 4339.7144
$ws.Range("J107").Value = 1997.5
$ws.Range("L107").Value = 1997.5
$ws.Range("N107").Value = -5837.5
$ws.Range("H132").Value = 977.5
$ws.Range("I132").Value = 977.5
$ws.Range("K132").Value = 2932.5
$ws.Range("M132").Value = -402.5
$ws.Range("H137").Value = 1106
$ws.Range("I137").Value = 1106
$ws.Range("K137").Value = 3318
$ws.Range("M137").Value = -768
$ws.Range("H141").Value = 1821.6666
$ws.Range("I141").Value = 1821.6666
$ws.Range("K141").Value = 5464.9998
$ws.Range("M141").Value = -284.9997999999996

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3798
$ws.Range("I2").Value = 3798
$ws.Range("K2").Value = 3798
$ws.Range("M2").Value = -3685
$ws.Range("H61").Value = 2366.3333
$ws.Range("I61").Value = 2366.3333
$ws.Range("K61").Value = 2366.3333
$ws.Range("M61").Value = -2154.3333
$ws.Range("H116").Value = 3798
$ws.Range("I116").Value = 3798
$ws.Range("K116").Value = 3798
$ws.Range("M116").Value = -1504
$ws.Range("H136").Value = 2366.3333
$ws.Range("I136").Value = 2366.3333
$ws.Range("K136").Value = 7098.999899999999
$ws.Range("M136").Value = -4548.999899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3798
$ws.Range("I3").Value = 3798
$ws.Range("K3").Value = 3798
$ws.Range("M3").Value = -3684

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1812.2632
$ws.Range("I31").Value = 1927.3334
$ws.Range("J31").Value = 1708.7
$ws.Range("K31").Value = 1927.3334
$ws.Range("L31").Value = 1708.7
$ws.Range("M31").Value = -1632.3334
$ws.Range("N31").Value = -2298.7
$ws.Range("H34").Value = 1812.2632
$ws.Range("I34").Value = 1927.3334
$ws.Range("J34").Value = 1708.7
$ws.Range("K34").Value = 1927.3334
$ws.Range("L34").Value = 1708.7
$ws.Range("M34").Value = -1725.3334
$ws.Range("N34").Value = -2112.7
$ws.Range("H105").Value = 4227
$ws.Range("I105").Value = 4658.6
$ws.Range("J105").Value = 3795.4
$ws.Range("K105").Value = 4658.6
$ws.Range("L105").Value = 3795.4
$ws.Range("M105").Value = -2911.6
$ws.Range("N105").Value = -7289.4
$ws.Range("H132").Value = 2738.3845
$ws.Range("I132").Value = 1951.75
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 5855.25
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -3325.25
$ws.Range("N132").Value = -17051
$ws.Range("H134").Value = 11248
$ws.Range("I134").Value = 11659.5
$ws.Range("K134").Value = 34978.5
$ws.Range("M134").Value = -32443.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 20
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = $null
$ws.Range("H18").Value = 2422.5
$ws.Range("I18").Value = 2230
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 6690
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = -6521
$ws.Range("N18").Value = -9338
$ws.Range("H141").Value = 7826
$ws.Range("I141").Value = 7826
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 23478
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = -18298

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = $null
$ws.Range("N47").Value = 0
$ws.Range("H126").Value = 6399.8
$ws.Range("I126").Value = 6399.8
$ws.Range("K126").Value = 19199.4
$ws.Range("M126").Value = -16729.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5730.778
$ws.Range("J22").Value = 7630
$ws.Range("L22").Value = 7630
$ws.Range("N22").Value = -8220
$ws.Range("H27").Value = 5730.778
$ws.Range("J27").Value = 7630
$ws.Range("L27").Value = 7630
$ws.Range("N27").Value = -7844
$ws.Range("H68").Value = 2762.818
$ws.Range("I68").Value = 2839.1
$ws.Range("K68").Value = 2839.1
$ws.Range("M68").Value = -2090.1
$ws.Range("H71").Value = 2762.818
$ws.Range("I71").Value = 2839.1
$ws.Range("K71").Value = 14195.5
$ws.Range("M71").Value = -10451.5
$ws.Range("H88").Value = 46666.332
$ws.Range("J88").Value = 46666.332
$ws.Range("L88").Value = 46666.332
$ws.Range("N88").Value = -47522.332
$ws.Range("H91").Value = 46666.332
$ws.Range("J91").Value = 46666.332
$ws.Range("L91").Value = 46666.332
$ws.Range("N91").Value = -49630.332
$ws.Range("H119").Value = 56105
$ws.Range("J119").Value = 56105
$ws.Range("L119").Value = 56105
$ws.Range("N119").Value = -65781
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = $null
$ws.Range("N120").Value = 0
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494
$ws.Range("H136").Value = 4326.6665
$ws.Range("I136").Value = 4326.6665
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12979.9995
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10429.9995
$ws.Range("N136").Value = $null

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1519.2307
$ws.Range("I107").Value = 1450.2
$ws.Range("J107").Value = 1562.375
$ws.Range("K107").Value = 4350.6
$ws.Range("L107").Value = 4687.125
$ws.Range("M107").Value = -2430.6
$ws.Range("N107").Value = -8527.125
$ws.Range("H132").Value = 2015.4286
$ws.Range("I132").Value = 1201.8182
$ws.Range("J132").Value = 4998.6665
$ws.Range("K132").Value = 3605.4546
$ws.Range("L132").Value = 14995.9995
$ws.Range("M132").Value = -1075.4546
$ws.Range("N132").Value = -20055.9995
$ws.Range("H136").Value = 6704.1
$ws.Range("I136").Value = 5434.5713
$ws.Range("K136").Value = 16303.7139
$ws.Range("M136").Value = -13753.7139
